{"js": "// The document contains three occurrences of an \"<id>...</id>\" marker,\n// each currently split across three runs, e.g. \"<id>\" + \"p090r_a1\" + \"</id>\".\n// The edit collapses each occurrence into a single run whose text drops the\n// \"a\" in front of the trailing digit: \"<id>p090r_a1</id>\" -> \"<id>p090r_1</id>\".\nconst replacements = [\n  [\"<id>p090r_a1</id>\", \"<id>p090r_1</id>\"],\n  [\"<id>p090r_a2</id>\", \"<id>p090r_2</id>\"],\n  [\"<id>p090r_a3</id>\", \"<id>p090r_3</id>\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# The document contains three \"<id>...</id>\" markers, each currently split\n# across three runs, e.g. \"<id>\" + \"p090r_a1\" + \"</id>\". Collapse each into\n# a single run whose text drops the \"a\" before the trailing digit:\n# \"<id>p090r_a1</id>\" -> \"<id>p090r_1</id>\".\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @{ Old = \"<id>p090r_a1</id>\"; New = \"<id>p090r_1</id>\" },\n    @{ Old = \"<id>p090r_a2</id>\"; New = \"<id>p090r_2</id>\" },\n    @{ Old = \"<id>p090r_a3</id>\"; New = \"<id>p090r_3</id>\" }\n)\n\nforeach ($pair in $pairs) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $pair.New\n    $find.Forward = $true\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute(\n        $pair.Old,      # FindText\n        $true,          # MatchCase\n        $false,         # MatchWholeWord\n        $false,         # MatchWildcards\n        $false,         # MatchSoundsLike\n        $false,         # MatchAllWordForms\n        $true,          # Forward\n        1,              # Wrap (wdFindContinue)\n        $false,         # Format\n        $pair.New,      # ReplaceWith\n        2               # Replace (wdReplaceAll)\n    ) | Out-Null\n}\n"}
